# Applies the NEAT_schema.xlsx restructuring:
#  - swap the Metric sheet's metric__name / metric__type headers
#  - replace ClassifierParams' six columns with two new summary columns
#    (classifierParams__sklearn_params, classifierParams__tf_keras_params)
#    and drop the old "optimizer" data validation
#  - insert a new "SkLearnParams" sheet right before "Target"
#  - insert a new "TFKerasParams" sheet right before "TrainValidData"
#    (carrying the old layers/loss/metrics/optimizer/fit_config headers
#    plus the optimizer dropdown validation)
#  - keep "Upload" the active/selected sheet, as it was originally

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metric: swap A1 / B1 header text.
# ---------------------------------------------------------------------
$metric = $wb.Worksheets.Item("Metric")
$metric.Range("A1").Value = "metric__type"
$metric.Range("B1").Value = "metric__name"

# ---------------------------------------------------------------------
# 2. ClassifierParams: drop the old optimizer validation + columns
#    C:F, rewrite A1:B1 to the two new summary headers.
# ---------------------------------------------------------------------
$classifierParams = $wb.Worksheets.Item("ClassifierParams")
$classifierParams.Range("E2:E1048576").Validation.Delete()
$classifierParams.Range("C1:F1").ClearContents()
$classifierParams.Range("A1").Value = "classifierParams__sklearn_params"
$classifierParams.Range("B1").Value = "classifierParams__tf_keras_params"

# ---------------------------------------------------------------------
# 3. New sheet "SkLearnParams", inserted right before "Target".
# ---------------------------------------------------------------------
$target = $wb.Worksheets.Item("Target")
$skLearnParams = $wb.Worksheets.Add($target)
$skLearnParams.Name = "SkLearnParams"
$skLearnParams.Range("A1").Value = "skLearnParams__random_state"
$skLearnParams.Range("B1").Value = "skLearnParams__max_iter"

# ---------------------------------------------------------------------
# 4. New sheet "TFKerasParams", inserted right before "TrainValidData",
#    re-using the old ClassifierParams layout + optimizer validation.
# ---------------------------------------------------------------------
$trainValidData = $wb.Worksheets.Item("TrainValidData")
$tfKerasParams = $wb.Worksheets.Add($trainValidData)
$tfKerasParams.Name = "TFKerasParams"
$tfKerasParams.Range("A1").Value = "tFKerasParams__layers_config"
$tfKerasParams.Range("B1").Value = "tFKerasParams__loss"
$tfKerasParams.Range("C1").Value = "tFKerasParams__metrics_config"
$tfKerasParams.Range("D1").Value = "tFKerasParams__optimizer"
$tfKerasParams.Range("E1").Value = "tFKerasParams__fit_config"
$tfKerasParams.Range("D2:D1048576").Validation.Add(3, 1, 1, '"adagrad,adam,adamax,nadam,sgd"')

# ---------------------------------------------------------------------
# 5. Restore the originally-selected/active sheet ("Upload" was the
#    last tab and stays the last tab after the two inserts above).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Upload").Activate()
